$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellFromTemplate {
    param($ws, $srcAddr, $dstAddr, $value)
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range($dstAddr).Value = $value
}

# --- Step 1: give the new last row (31) the distinctive "final row" border
#     that row 25 currently has, copying per-cell so no stray cells appear ---
Set-CellFromTemplate $ws "A25" "A31" 45572.63892673611
Set-CellFromTemplate $ws "B25" "B31" 'BetaLife'
Set-CellFromTemplate $ws "C25" "C31" 'Yes'
Set-CellFromTemplate $ws "D25" "D31" 'Yes'
Set-CellFromTemplate $ws "E25" "E31" 'Yes'
Set-CellFromTemplate $ws "F25" "F31" 'Yes'
Set-CellFromTemplate $ws "G25" "G31" 'Yes'
Set-CellFromTemplate $ws "H25" "H31" 'Yes'
Set-CellFromTemplate $ws "I25" "I31" 'At least 1 deployment is running in a live/production environment (Scaling)'
Set-CellFromTemplate $ws "J25" "J31" 'Active deployments or customers in 1 to 3 countries (Scaling)'
Set-CellFromTemplate $ws "K25" "K31" 'Less than 100 unique daily users (Pilot)'
Set-CellFromTemplate $ws "L25" "L31" 'Between 1000 and 10000 (Scaling)'
Set-CellFromTemplate $ws "M25" "M31" 'Greater than 150,000 USD (Mature)'
Set-CellFromTemplate $ws "N25" "N31" 'Between 50,000 and 250,000 USD (Scaling)'

# --- Step 2: row 25 is no longer the last row, so give it the normal
#     alternating-row style (template: row 23) while preserving its values ---
Set-CellFromTemplate $ws "A23" "A25" 45568.454662511576
Set-CellFromTemplate $ws "B23" "B25" 'Medpharma'
Set-CellFromTemplate $ws "C23" "C25" 'Yes'
Set-CellFromTemplate $ws "D23" "D25" 'Yes'
Set-CellFromTemplate $ws "E23" "E25" 'Yes'
Set-CellFromTemplate $ws "F23" "F25" 'Yes'
Set-CellFromTemplate $ws "G23" "G25" 'Yes'
Set-CellFromTemplate $ws "H23" "H25" 'Yes'
Set-CellFromTemplate $ws "I23" "I25" 'At least 1 deployment is running in a live/production environment (Scaling)'
Set-CellFromTemplate $ws "J23" "J25" 'Active deployments or customers in 1 to 3 countries (Scaling)'
Set-CellFromTemplate $ws "K23" "K25" 'Less than 100 unique daily users (Pilot)'
Set-CellFromTemplate $ws "L23" "L25" 'Less than 1000 (Pilot)'
Set-CellFromTemplate $ws "M23" "M25" 'Greater than 150,000 USD (Mature)'
Set-CellFromTemplate $ws "N23" "N25" 'Between 50,000 and 250,000 USD (Scaling)'
Set-CellFromTemplate $ws "O23" "O25" 'No'
Set-CellFromTemplate $ws "P23" "P25" 'Yes'
Set-CellFromTemplate $ws "Q23" "Q25" 'No'
Set-CellFromTemplate $ws "R23" "R25" 'Yes'
Set-CellFromTemplate $ws "S23" "S25" 'Unknown'
Set-CellFromTemplate $ws "T23" "T25" 'Yes'
Set-CellFromTemplate $ws "U23" "U25" 'Unknown'
Set-CellFromTemplate $ws "V23" "V25" 'Yes'
Set-CellFromTemplate $ws "W23" "W25" '3 or fewer developers currently working on development (medium)'
Set-CellFromTemplate $ws "X23" "X25" 'More than 2 releases in past 12 months (high)'
Set-CellFromTemplate $ws "Y23" "Y25" 'No'
Set-CellFromTemplate $ws "Z23" "Z25" 'No'
Set-CellFromTemplate $ws "AA23" "AA25" 'Yes'

# --- Step 3: new response rows, copying cell-by-cell format from the nearest
#     same-parity existing row so the banding / borders keep matching ---
# Row 26 (template row 24)
Set-CellFromTemplate $ws "A24" "A26" 45569.53350535879
Set-CellFromTemplate $ws "B24" "B26" 'Alô Vida+'
Set-CellFromTemplate $ws "C24" "C26" 'Yes'
Set-CellFromTemplate $ws "D24" "D26" 'Yes'
Set-CellFromTemplate $ws "E24" "E26" 'Yes'
Set-CellFromTemplate $ws "F24" "F26" 'Yes'
Set-CellFromTemplate $ws "G24" "G26" 'Yes'
Set-CellFromTemplate $ws "H24" "H26" 'Yes'
Set-CellFromTemplate $ws "I24" "I26" 'At least 1 deployment is running in a live/production environment (Scaling)'
Set-CellFromTemplate $ws "J24" "J26" 'Active deployments or customers in 1 to 3 countries (Scaling)'
Set-CellFromTemplate $ws "K24" "K26" 'Between 100 and 1000 unique daily users (Scaling)'
Set-CellFromTemplate $ws "L24" "L26" 'Between 1000 and 10000 (Scaling)'
Set-CellFromTemplate $ws "M24" "M26" 'Greater than 150,000 USD (Mature)'
Set-CellFromTemplate $ws "N24" "N26" 'Over 250,000 USD (Mature)'
Set-CellFromTemplate $ws "O24" "O26" 'No'
Set-CellFromTemplate $ws "P24" "P26" 'No'
Set-CellFromTemplate $ws "Q24" "Q26" 'Yes'
Set-CellFromTemplate $ws "R24" "R26" 'Yes'
Set-CellFromTemplate $ws "S24" "S26" 'Yes'
Set-CellFromTemplate $ws "T24" "T26" 'Yes'
Set-CellFromTemplate $ws "U24" "U26" 'Yes'
Set-CellFromTemplate $ws "V24" "V26" 'Yes'
Set-CellFromTemplate $ws "W24" "W26" '3 or fewer developers currently working on development (medium)'
Set-CellFromTemplate $ws "X24" "X26" 'More than 2 releases in past 12 months (high)'
Set-CellFromTemplate $ws "Y24" "Y26" 'No'
Set-CellFromTemplate $ws "Z24" "Z26" 'Yes'
Set-CellFromTemplate $ws "AA24" "AA26" 'Yes'

# Row 27 (template row 23)
Set-CellFromTemplate $ws "A23" "A27" 45569.69487824074
Set-CellFromTemplate $ws "B23" "B27" 'MedTrack'
Set-CellFromTemplate $ws "C23" "C27" 'Yes'
Set-CellFromTemplate $ws "D23" "D27" 'Yes'
Set-CellFromTemplate $ws "E23" "E27" 'Yes'
Set-CellFromTemplate $ws "F23" "F27" 'Yes'
Set-CellFromTemplate $ws "G23" "G27" 'Yes'
Set-CellFromTemplate $ws "H23" "H27" 'Yes'
Set-CellFromTemplate $ws "I23" "I27" 'At least 1 deployment is running in a live/production environment (Scaling)'
Set-CellFromTemplate $ws "J23" "J27" 'Active deployments or customers in 1 to 3 countries (Scaling)'
Set-CellFromTemplate $ws "K23" "K27" 'Between 100 and 1000 unique daily users (Scaling)'
Set-CellFromTemplate $ws "L23" "L27" 'Between 1000 and 10000 (Scaling)'
Set-CellFromTemplate $ws "M23" "M27" 'Between 25,000 and 150,000 USD (Scaling)'
Set-CellFromTemplate $ws "N23" "N27" 'Over 250,000 USD (Mature)'
Set-CellFromTemplate $ws "O23" "O27" 'Yes'
Set-CellFromTemplate $ws "Q23" "Q27" 'Yes'
Set-CellFromTemplate $ws "R23" "R27" 'Yes'
Set-CellFromTemplate $ws "S23" "S27" 'Unknown'
Set-CellFromTemplate $ws "T23" "T27" 'Yes'
Set-CellFromTemplate $ws "U23" "U27" 'Yes'
Set-CellFromTemplate $ws "V23" "V27" 'Yes'
Set-CellFromTemplate $ws "W23" "W27" 'More than 3 developers currently working (high)'
Set-CellFromTemplate $ws "X23" "X27" 'More than 2 releases in past 12 months (high)'
Set-CellFromTemplate $ws "Y23" "Y27" 'No'
Set-CellFromTemplate $ws "Z23" "Z27" 'No'
Set-CellFromTemplate $ws "AA23" "AA27" 'Yes'

# Row 28 (template row 24)
Set-CellFromTemplate $ws "A24" "A28" 45572.434264166666
Set-CellFromTemplate $ws "B24" "B28" 'Spes 360'
Set-CellFromTemplate $ws "C24" "C28" 'Yes'
Set-CellFromTemplate $ws "D24" "D28" 'Yes'
Set-CellFromTemplate $ws "E24" "E28" 'Yes'
Set-CellFromTemplate $ws "F24" "F28" 'Yes'
Set-CellFromTemplate $ws "G24" "G28" 'Yes'
Set-CellFromTemplate $ws "H24" "H28" 'Yes'
Set-CellFromTemplate $ws "I24" "I28" 'At least 1 deployment is running in a live/production environment (Scaling)'
Set-CellFromTemplate $ws "J24" "J28" 'Active deployments or customers in 1 to 3 countries (Scaling)'
Set-CellFromTemplate $ws "K24" "K28" 'Between 100 and 1000 unique daily users (Scaling)'
Set-CellFromTemplate $ws "L24" "L28" 'Less than 1000 (Pilot)'
Set-CellFromTemplate $ws "M24" "M28" 'Between 25,000 and 150,000 USD (Scaling)'
Set-CellFromTemplate $ws "N24" "N28" 'Over 250,000 USD (Mature)'
Set-CellFromTemplate $ws "O24" "O28" 'No'
Set-CellFromTemplate $ws "P24" "P28" 'No'
Set-CellFromTemplate $ws "Q24" "Q28" 'No'
Set-CellFromTemplate $ws "R24" "R28" 'No'
Set-CellFromTemplate $ws "S24" "S28" 'Unknown'
Set-CellFromTemplate $ws "T24" "T28" 'Yes'
Set-CellFromTemplate $ws "U24" "U28" 'Yes'
Set-CellFromTemplate $ws "V24" "V28" 'Yes'
Set-CellFromTemplate $ws "W24" "W28" '3 or fewer developers currently working on development (medium)'
Set-CellFromTemplate $ws "X24" "X28" 'More than 2 releases in past 12 months (high)'
Set-CellFromTemplate $ws "Y24" "Y28" 'No'
Set-CellFromTemplate $ws "Z24" "Z28" 'No'
Set-CellFromTemplate $ws "AA24" "AA28" 'No'

# Row 29 (template row 23)
Set-CellFromTemplate $ws "A23" "A29" 45572.626628275466
Set-CellFromTemplate $ws "B23" "B29" 'Aviro Pocket Clinic'
Set-CellFromTemplate $ws "C23" "C29" 'Yes'
Set-CellFromTemplate $ws "D23" "D29" 'Yes'
Set-CellFromTemplate $ws "E23" "E29" 'Yes'
Set-CellFromTemplate $ws "F23" "F29" 'Yes'
Set-CellFromTemplate $ws "G23" "G29" 'Yes'
Set-CellFromTemplate $ws "H23" "H29" 'Yes'
Set-CellFromTemplate $ws "I23" "I29" 'At least 1 deployment is running in a live/production environment (Scaling)'
Set-CellFromTemplate $ws "J23" "J29" 'Active deployments or customers in 1 to 3 countries (Scaling)'
Set-CellFromTemplate $ws "K23" "K29" 'Less than 100 unique daily users (Pilot)'
Set-CellFromTemplate $ws "L23" "L29" 'Less than 1000 (Pilot)'
Set-CellFromTemplate $ws "M23" "M29" 'Greater than 150,000 USD (Mature)'
Set-CellFromTemplate $ws "N23" "N29" 'Over 250,000 USD (Mature)'
Set-CellFromTemplate $ws "O23" "O29" 'Yes'
Set-CellFromTemplate $ws "P23" "P29" 'Yes'
Set-CellFromTemplate $ws "Q23" "Q29" 'Yes'
Set-CellFromTemplate $ws "R23" "R29" 'Yes'
Set-CellFromTemplate $ws "S23" "S29" 'Unknown'
Set-CellFromTemplate $ws "T23" "T29" 'Yes'
Set-CellFromTemplate $ws "U23" "U29" 'Unknown'
Set-CellFromTemplate $ws "V23" "V29" 'Unknown'
Set-CellFromTemplate $ws "W23" "W29" 'No current active development (low)'
Set-CellFromTemplate $ws "X23" "X29" '1-2 releases in the past 12 months (medium)'
Set-CellFromTemplate $ws "Y23" "Y29" 'Yes'
Set-CellFromTemplate $ws "Z23" "Z29" 'Yes'
Set-CellFromTemplate $ws "AA23" "AA29" 'Yes'

# Row 30 (template row 24)
Set-CellFromTemplate $ws "A24" "A30" 45572.63018476852
Set-CellFromTemplate $ws "B24" "B30" 'XanaHealth '
Set-CellFromTemplate $ws "C24" "C30" 'Yes'
Set-CellFromTemplate $ws "D24" "D30" 'Yes'
Set-CellFromTemplate $ws "E24" "E30" 'Yes'
Set-CellFromTemplate $ws "F24" "F30" 'Yes'
Set-CellFromTemplate $ws "G24" "G30" 'Yes'
Set-CellFromTemplate $ws "H24" "H30" 'Yes'
Set-CellFromTemplate $ws "I24" "I30" 'At least 1 deployment is running in a live/production environment (Scaling)'
Set-CellFromTemplate $ws "J24" "J30" 'Active deployments or customers in 1 to 3 countries (Scaling)'
Set-CellFromTemplate $ws "K24" "K30" 'Between 100 and 1000 unique daily users (Scaling)'
Set-CellFromTemplate $ws "L24" "L30" 'Between 1000 and 10000 (Scaling)'
Set-CellFromTemplate $ws "M24" "M30" 'Between 25,000 and 150,000 USD (Scaling)'
Set-CellFromTemplate $ws "N24" "N30" 'Between 50,000 and 250,000 USD (Scaling)'
Set-CellFromTemplate $ws "O24" "O30" 'Yes'
Set-CellFromTemplate $ws "P24" "P30" 'Yes'
Set-CellFromTemplate $ws "Q24" "Q30" 'Yes'
Set-CellFromTemplate $ws "R24" "R30" 'Yes'
Set-CellFromTemplate $ws "S24" "S30" 'Unknown'
Set-CellFromTemplate $ws "T24" "T30" 'Yes'
Set-CellFromTemplate $ws "U24" "U30" 'Unknown'
Set-CellFromTemplate $ws "V24" "V30" 'Unknown'
Set-CellFromTemplate $ws "W24" "W30" 'More than 3 developers currently working (high)'
Set-CellFromTemplate $ws "X24" "X30" 'More than 2 releases in past 12 months (high)'
Set-CellFromTemplate $ws "Y24" "Y30" 'Yes'
Set-CellFromTemplate $ws "Z24" "Z30" 'Yes'
Set-CellFromTemplate $ws "AA24" "AA30" 'Yes'

# --- Step 4: expand the "Form_Responses1" table to include the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:AA31"))
